$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Construction Safety and Compliance"
$ws.Range("B12").Value = "abc"

$ws.Range("A13").Value = "Sustainable Building and Construction"
$ws.Range("B13").Value = "xyz"

$ws.Range("A14").Value = "Sustainable Building and Construction"
$ws.Range("B14").Value = "abc"
